$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-08-31"

# Update the August header label to reflect the new "through" date
$ws.Range("A9").Value = "August (through 08-31)"

# Update August row (row 9) values
$ws.Range("B9").Value = 32
$ws.Range("C9").Value = 79
$ws.Range("F9").Value = 45
$ws.Range("G9").Value = 163
$ws.Range("H9").Value = 160
$ws.Range("I9").Value = 167

# Update Total row (row 10) values
$ws.Range("B10").Value = 194
$ws.Range("C10").Value = 381
$ws.Range("F10").Value = 349
$ws.Range("G10").Value = 784
$ws.Range("H10").Value = 1070
$ws.Range("I10").Value = 1138
